# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the crypto symbol
# list with the latest scraped quotes (GitHub Actions run of Sat Jan 7
# 23:45:50 UTC 2023). Values are stored as text in the sheet, so each target
# cell is forced to the "@" (text) number format before the new value is
# written — this avoids Excel's automatic number/percentage re-typing of the
# numeric-looking strings, keeping e.g. trailing zeros and leading "0."
# exactly as scraped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.30"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.68%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.13"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.94%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.708"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.95%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06211"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.725"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8506"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.23%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9067"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.67%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1403"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.02%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04727"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-11.01%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07093"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.30%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03173"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.93%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09056"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.88%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001530"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.65%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006157"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.63%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005992"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.89%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.42%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.28%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.177"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.38%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.68%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.58%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.103"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.11%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04246"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.00%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001221"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.29%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004117"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.04%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.07%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03898"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.51%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.19%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004133"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.07%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.76%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01341"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-10.16%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005174"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.86%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.06%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.03503"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-35.80%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05783"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-56.23%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.06%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.06%"
